# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '73.913.46'
$ws.Range("E2").Value = '  +7.43%  '
$ws.Range("D3").Value = '2.619.96'
$ws.Range("E3").Value = '  +7.22%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '187.02'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +14.17%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '582.28'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.83%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +4.66%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.199'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +16.54%  '
$ws.Range("D10").Value = '2.619.08'
$ws.Range("E10").Value = '  +7.30%  '
$ws.Range("E11").Value = '  +1.29%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.358'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +7.58%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.69'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.89%  '
# Row 14/15: ShibaInu and WrappedBTC swapped ranking positions this refresh.
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.0000189'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +5.38%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '73.881.12'
$ws.Range("E15").Value = '  +7.54%  '
$ws.Range("E16").Value = '  +7.37%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '26.44'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +12.79%  '
$ws.Range("D18").Value = '2.631.87'
$ws.Range("E18").Value = '  +7.74%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '9.11'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +29.81%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '11.78'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +10.93%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '366.61'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +8.07%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.30'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +18.06%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.08'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +5.92%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.40%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '69.95'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +6.94%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '4.13'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +8.85%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.31'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +10.88%  '
$ws.Range("D28").Value = '2.754.21'
$ws.Range("E28").Value = '  +7.24%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").Value = '0.0₃0941'
$ws.Range("E30").Value = '  +14.05%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '523.12'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +20.67%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.39'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +15.09%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '7.67'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +6.54%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.75'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +9.04%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '162.57'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.05%  '
$ws.Range("E37").Value = '  +10.60%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '19.11'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +6.14%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '19.27'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("E40").Value = '  +0.06%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.92'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +12.14%  '
$ws.Range("E42").Value = '  +9.67%  '
$ws.Range("E43").Value = '  +8.05%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '161.60'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +24.20%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.38'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +14.05%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +8.98%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '38.88'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +3.56%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0854'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +18.59%  '
$ws.Range("E49").Value = '  +8.09%  '
$ws.Range("E50").Value = '  +7.77%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '20.82'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +22.53%  '
